# Refresh the crypto symbol list: updated prices/volume labels and bumped
# the "Hora" (hour) column from 2 -> 3 for every data row, matching the
# GitHub Actions scrape commit on Sat Dec 17 03:07:53 UTC 2022.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) and G (Hora) store plain numeric-looking strings as TEXT
# in this workbook, so a leading apostrophe is used to keep Excel from
# reinterpreting them as numbers when the value is written back.

$ws.Range("D2").Value = "'229.14"
$ws.Range("G2").Value = "'3"

$ws.Range("D3").Value = "'22.50"
$ws.Range("G3").Value = "'3"

$ws.Range("D4").Value = "'5.343"
$ws.Range("G4").Value = "'3"

$ws.Range("D5").Value = "'0.05528"
$ws.Range("G5").Value = "'3"

$ws.Range("D6").Value = "'3.390"
$ws.Range("G6").Value = "'3"

$ws.Range("D7").Value = "'6.469"
$ws.Range("G7").Value = "'3"

$ws.Range("D8").Value = "'0.7815"
$ws.Range("G8").Value = "'3"

$ws.Range("D9").Value = "'1.044"
$ws.Range("G9").Value = "'3"

$ws.Range("D10").Value = "'0.1382"
$ws.Range("G10").Value = "'3"

$ws.Range("D11").Value = "'0.07446"
$ws.Range("G11").Value = "'3"

$ws.Range("G12").Value = "'3"

$ws.Range("D13").Value = "'0.02945"
$ws.Range("G13").Value = "'3"

$ws.Range("G14").Value = "'3"

$ws.Range("D15").Value = "'0.001660"
$ws.Range("G15").Value = "'3"

$ws.Range("D16").Value = "'3.255"
$ws.Range("G16").Value = "'3"

$ws.Range("D17").Value = "'0.04777"
$ws.Range("G17").Value = "'3"

$ws.Range("G18").Value = "'3"

$ws.Range("D19").Value = "'0.006198"
$ws.Range("G19").Value = "'3"

$ws.Range("D20").Value = "'0.005236"
$ws.Range("G20").Value = "'3"

$ws.Range("G21").Value = "'3"

$ws.Range("G22").Value = "'3"

$ws.Range("D23").Value = "'3.891"
$ws.Range("G23").Value = "'3"

$ws.Range("G24").Value = "'3"

$ws.Range("G25").Value = "'3"

$ws.Range("D26").Value = "'0.1281"
$ws.Range("G26").Value = "'3"

$ws.Range("E27").Value = "26UpBotsUBXTBestin24h"
$ws.Range("G27").Value = "'3"

$ws.Range("G28").Value = "'3"
$ws.Range("G29").Value = "'3"
$ws.Range("G30").Value = "'3"
$ws.Range("G31").Value = "'3"
$ws.Range("G32").Value = "'3"
$ws.Range("G33").Value = "'3"
$ws.Range("G34").Value = "'3"
$ws.Range("G35").Value = "'3"
$ws.Range("G36").Value = "'3"
$ws.Range("G37").Value = "'3"
$ws.Range("G38").Value = "'3"
$ws.Range("G39").Value = "'3"

$ws.Range("D40").Value = "'0.03949"
$ws.Range("G40").Value = "'3"

$ws.Range("D41").Value = "'0.007138"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("G41").Value = "'3"

$ws.Range("D42").Value = "'0.1035"
$ws.Range("G42").Value = "'3"

$ws.Range("D43").Value = "'0.003170"
$ws.Range("G43").Value = "'3"

$ws.Range("D44").Value = "'0.009260"
$ws.Range("G44").Value = "'3"

$ws.Range("D45").Value = "'0.00005425"
$ws.Range("G45").Value = "'3"

$ws.Range("G46").Value = "'3"
$ws.Range("G47").Value = "'3"

$ws.Range("D48").Value = "'0.08837"
$ws.Range("G48").Value = "'3"

$ws.Range("G49").Value = "'3"
$ws.Range("G50").Value = "'3"
$ws.Range("G51").Value = "'3"
